$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Version History")

# Row 2: update author / updated-section text, and replace the calculated
# Date column with a plain (no longer formula-driven) value.
$ws.Range("B2").Value = "Omar Sherif "
$ws.Range("C2").Value = "Registerations test cases update after review "
$ws.Range("D2").Value = 45768

# Rows 3-5: the older version-history entries are removed. Clear their
# contents and let the row heights return to the sheet default (matches
# the blank filler rows 6-9 already on the sheet).
$ws.Range("A3:D5").ClearContents()
$ws.Rows("3:5").EntireRow.AutoFit()

# Update the selected cell to match the saved view state.
$ws.Range("D17").Select()
